$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header cells (I0, IF) in row 1, copying the header
# formatting (style) from the existing H1 header cell.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "IF"

# Data for the new columns I (I0) and J (IF), one triple per data row:
# (row number, I value, J value)
$data = @(
    @(2, 7, 7),
    @(3, 8, 9),
    @(4, 5, 5),
    @(5, 3, 5),
    @(6, 9, 9),
    @(7, 7, 8),
    @(8, 7, 8),
    @(9, 4, 6),
    @(10, 7, 8),
    @(11, 7, 7),
    @(12, 7, 8),
    @(13, 6, 8),
    @(14, 6, 8),
    @(15, 2, 7),
    @(16, 1, 2),
    @(17, 1, 2),
    @(18, 1, 6),
    @(19, 1, 7),
    @(20, 1, 2),
    @(21, 1, 6),
    @(22, 1, 5),
    @(23, 1, 6),
    @(24, 1, 6),
    @(25, 1, 5),
    @(26, 1, 5),
    @(27, 1, 5),
    @(28, 5, 7),
    @(29, 3, 5),
    @(30, 8, 8),
    @(31, 5, 6),
    @(32, 4, 6),
    @(33, 9, 9),
    @(34, 7, 7),
    @(35, 7, 7),
    @(36, 8, 8),
    @(37, 5, 6),
    @(38, 7, 7),
    @(39, 7, 7),
    @(40, 8, 9),
    @(41, 7, 8),
    @(42, 5, 6),
    @(43, 6, 7),
    @(44, 6, 7),
    @(45, 6, 7),
    @(46, 1, 3),
    @(47, 1, 4),
    @(48, 1, 6),
    @(49, 1, 5),
    @(50, 1, 3),
    @(51, 1, 3),
    @(52, 1, 2)
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 9).Value = $item[1]
    $ws.Cells.Item($r, 10).Value = $item[2]
}
